$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Cost" column header and its values
$ws.Range("F1").Value = "Cost"
$ws.Range("F2").Value = 50
$ws.Range("F3").Value = 300
$ws.Range("F4").Value = 500

# Match the new column sizing seen after the edit (Excel auto-fit side effect
# of the columns now containing shorter/longer text once Cost was added)
$ws.Columns.Item(1).ColumnWidth = 5.785714285714286
$ws.Columns.Item(2).ColumnWidth = 2.4107142857142856
$ws.Columns.Item(3).ColumnWidth = 9.910714285714286
$ws.Columns.Item(5).ColumnWidth = 19.035714285714285

# Move the active selection as recorded after the edit
$ws.Range("I11").Select()
